$wb = $excel.ActiveWorkbook

$wsAbertos  = $wb.Worksheets.Item("ContravalesAbertos")
$wsBaixados = $wb.Worksheets.Item("ContravalesBaixados")

# Contravales sendo baixados (settled) nesta operacao, com a data/hora da baixa.
$baixas = @(
    @{ Numero = 1; DataHora = "2025-03-25 15:58:44" },
    @{ Numero = 2; DataHora = "2025-03-25 15:59:26" },
    @{ Numero = 4; DataHora = "2025-03-25 15:59:57" }
)

foreach ($baixa in $baixas) {
    # Encontra a linha do contravale em aberto e remove (shift up).
    $lastRow = $wsAbertos.Cells.Item($wsAbertos.Rows.Count, 1).End(-4162).Row
    for ($r = $lastRow; $r -ge 2; $r--) {
        $valor = $wsAbertos.Cells.Item($r, 1).Value2
        if ($valor -eq $baixa.Numero) {
            $wsAbertos.Rows.Item($r).Delete()
            break
        }
    }

    # Adiciona o contravale baixado na planilha de baixados.
    $novaLinha = $wsBaixados.Cells.Item($wsBaixados.Rows.Count, 1).End(-4162).Row + 1
    $wsBaixados.Cells.Item($novaLinha, 1).Value = $baixa.Numero
    $wsBaixados.Cells.Item($novaLinha, 2).Value = $baixa.DataHora
}
